# Updates cryptos list data per diff (price + volume columns, and two row re-orderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "30.078.62"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +3.03%  "

# Row 3
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.868.38"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +2.33%  "

# Row 4
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9987"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "246.43"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "

# Row 6
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6438"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +4.05%  "

# Row 7
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9991"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3008"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +3.94%  "

# Row 9
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07516"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +2.31%  "

# Row 10
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "24.47"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +6.39%  "

# Row 11
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07691"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.871.86"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +2.69%  "

# Row 13
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "5.082"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +2.55%  "

# Row 14
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6936"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +4.67%  "

# Row 15
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "84.50"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +2.90%  "

# Row 16
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000009495"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +6.31%  "

# Row 17
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "6.122"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +4.83%  "

# Row 18
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "30.061.62"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +3.10%  "

# Row 19
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "2.120.32"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +2.63%  "

# Row 20
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "241.79"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "

# Row 21
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "12.74"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +2.51%  "

# Row 22
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9994"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "7.478"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +4.08%  "

# Row 24
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "1.000"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "159.74"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "

# Row 26
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1432"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "

# Row 27
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "8.625"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +2.36%  "

# Row 29
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06175"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +11.06%  "

# Row 30
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "1.510"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +1.84%  "

# Row 31
$ws.Range("E31").Value = "  +5.78%  "

# Row 32
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "4.181"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +2.01%  "

# Row 33
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "4.146"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "

# Row 34
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "1.880"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +3.29%  "

# Row 35
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.169"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +3.36%  "

# Row 36
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7367"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "

# Row 37
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "2.881"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "

# Row 39
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.01810"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +3.00%  "

# Row 40
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "1.226.03"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +1.05%  "

# Row 41
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9317"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +1.52%  "

# Row 42
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "6.299"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "

# Row 43
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "2.037.54"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +3.51%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "

# Row 45
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "102.50"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "66.80"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +3.29%  "

# Row 47
$ws.Range("E47").Value = "  -2.45%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "9.389"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +3.50%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5087"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "

# Row 50
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4115"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +2.61%  "

# Row 51
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1143"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +3.12%  "
